# Generate Report for Handoff
#
# A new handoff pass was run for the "27a0d64a-1b36-42a4-9fa8-38a783daef7e"
# file, which refreshes the "Latest Handoff Datetime" (column D, row 6) on
# the per-language status sheets. (Row 6's timestamp had previously been
# stuck duplicating row 7's value; this records the new, distinct handoff
# time for row 6 while row 7 keeps the datetime it already showed.)

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("D6").Value = "2016-03-09 10:37:04"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("D6").Value = "2016-03-09 10:37:09"
